$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 40.23443200000001
$ws.Cells.Item(2, 8).Value = 120.703296
$ws.Cells.Item(2, 9).Value = 0.5194057602668869
$ws.Cells.Item(2, 10).Value = 0.5194057602668869
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 40.23443200000001
$ws.Cells.Item(2, 14).Value = 120.703296
$ws.Cells.Item(2, 15).Value = 0.5194057602668869
$ws.Cells.Item(2, 16).Value = 0.5194057602668869
$ws.Cells.Item(2, 17).Value = 1618.809518362624
$ws.Cells.Item(2, 18).Value = 14569.28566526362
$ws.Cells.Item(2, 19).Value = 0.2697823437984228
$ws.Cells.Item(2, 20).Value = 0.2697823437984228
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 40.23443200000001
$ws.Cells.Item(3, 8).Value = 120.703296
$ws.Cells.Item(3, 9).Value = 0.5194057602668869
$ws.Cells.Item(3, 10).Value = 0.5194057602668869
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.585148666666667
$ws.Cells.Item(3, 14).Value = 7.755446
$ws.Cells.Item(3, 15).Value = 0.03337293561427507
$ws.Cells.Item(3, 16).Value = 0.03337293561427508
$ws.Cells.Item(3, 17).Value = 104.0119882388907
$ws.Cells.Item(3, 18).Value = 936.107894150016
$ws.Cells.Item(3, 19).Value = 0.01733409499507041
$ws.Cells.Item(3, 20).Value = 0.01733409499507042
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 40.23443200000001
$ws.Cells.Item(4, 8).Value = 120.703296
$ws.Cells.Item(4, 9).Value = 0.5194057602668869
$ws.Cells.Item(4, 10).Value = 0.5194057602668869
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.525608
$ws.Cells.Item(4, 14).Value = 1.576824
$ws.Cells.Item(4, 15).Value = 0.006785328120013172
$ws.Cells.Item(4, 16).Value = 0.006785328120013173
$ws.Cells.Item(4, 17).Value = 21.147539334656
$ws.Cells.Item(4, 18).Value = 190.327854011904
$ws.Cells.Item(4, 19).Value = 0.003524338510835728
$ws.Cells.Item(4, 20).Value = 0.003524338510835728
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 40.23443200000001
$ws.Cells.Item(5, 8).Value = 120.703296
$ws.Cells.Item(5, 9).Value = 0.5194057602668869
$ws.Cells.Item(5, 10).Value = 0.5194057602668869
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 34.11724066666667
$ws.Cells.Item(5, 14).Value = 102.351722
$ws.Cells.Item(5, 15).Value = 0.4404359759988248
$ws.Cells.Item(5, 16).Value = 0.4404359759988249
$ws.Cells.Item(5, 17).Value = 1372.687799630635
$ws.Cells.Item(5, 18).Value = 12354.19019667571
$ws.Cells.Item(5, 19).Value = 0.2287649829625579
$ws.Cells.Item(5, 20).Value = 0.228764982962558
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 2.585148666666667
$ws.Cells.Item(6, 8).Value = 7.755446
$ws.Cells.Item(6, 9).Value = 0.03337293561427507
$ws.Cells.Item(6, 10).Value = 0.03337293561427508
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 40.23443200000001
$ws.Cells.Item(6, 14).Value = 120.703296
$ws.Cells.Item(6, 15).Value = 0.5194057602668869
$ws.Cells.Item(6, 16).Value = 0.5194057602668869
$ws.Cells.Item(6, 17).Value = 104.0119882388907
$ws.Cells.Item(6, 18).Value = 936.107894150016
$ws.Cells.Item(6, 19).Value = 0.01733409499507041
$ws.Cells.Item(6, 20).Value = 0.01733409499507042
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 2.585148666666667
$ws.Cells.Item(7, 8).Value = 7.755446
$ws.Cells.Item(7, 9).Value = 0.03337293561427507
$ws.Cells.Item(7, 10).Value = 0.03337293561427508
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.585148666666667
$ws.Cells.Item(7, 14).Value = 7.755446
$ws.Cells.Item(7, 15).Value = 0.03337293561427507
$ws.Cells.Item(7, 16).Value = 0.03337293561427508
$ws.Cells.Item(7, 17).Value = 6.682993628768444
$ws.Cells.Item(7, 18).Value = 60.146942658916
$ws.Cells.Item(7, 19).Value = 0.00111375283151455
$ws.Cells.Item(7, 20).Value = 0.00111375283151455
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.585148666666667
$ws.Cells.Item(8, 8).Value = 7.755446
$ws.Cells.Item(8, 9).Value = 0.03337293561427507
$ws.Cells.Item(8, 10).Value = 0.03337293561427508
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.525608
$ws.Cells.Item(8, 14).Value = 1.576824
$ws.Cells.Item(8, 15).Value = 0.006785328120013172
$ws.Cells.Item(8, 16).Value = 0.006785328120013173
$ws.Cells.Item(8, 17).Value = 1.358774820389333
$ws.Cells.Item(8, 18).Value = 12.228973383504
$ws.Cells.Item(8, 19).Value = 0.0002264463184709297
$ws.Cells.Item(8, 20).Value = 0.0002264463184709298
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.585148666666667
$ws.Cells.Item(9, 8).Value = 7.755446
$ws.Cells.Item(9, 9).Value = 0.03337293561427507
$ws.Cells.Item(9, 10).Value = 0.03337293561427508
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 34.11724066666667
$ws.Cells.Item(9, 14).Value = 102.351722
$ws.Cells.Item(9, 15).Value = 0.4404359759988248
$ws.Cells.Item(9, 16).Value = 0.4404359759988249
$ws.Cells.Item(9, 17).Value = 88.19813921977911
$ws.Cells.Item(9, 18).Value = 793.7832529780121
$ws.Cells.Item(9, 19).Value = 0.01469864146921918
$ws.Cells.Item(9, 20).Value = 0.01469864146921919
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.525608
$ws.Cells.Item(10, 8).Value = 1.576824
$ws.Cells.Item(10, 9).Value = 0.006785328120013172
$ws.Cells.Item(10, 10).Value = 0.006785328120013173
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 40.23443200000001
$ws.Cells.Item(10, 14).Value = 120.703296
$ws.Cells.Item(10, 15).Value = 0.5194057602668869
$ws.Cells.Item(10, 16).Value = 0.5194057602668869
$ws.Cells.Item(10, 17).Value = 21.147539334656
$ws.Cells.Item(10, 18).Value = 190.327854011904
$ws.Cells.Item(10, 19).Value = 0.003524338510835728
$ws.Cells.Item(10, 20).Value = 0.003524338510835728
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.525608
$ws.Cells.Item(11, 8).Value = 1.576824
$ws.Cells.Item(11, 9).Value = 0.006785328120013172
$ws.Cells.Item(11, 10).Value = 0.006785328120013173
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 2.585148666666667
$ws.Cells.Item(11, 14).Value = 7.755446
$ws.Cells.Item(11, 15).Value = 0.03337293561427507
$ws.Cells.Item(11, 16).Value = 0.03337293561427508
$ws.Cells.Item(11, 17).Value = 1.358774820389333
$ws.Cells.Item(11, 18).Value = 12.228973383504
$ws.Cells.Item(11, 19).Value = 0.0002264463184709297
$ws.Cells.Item(11, 20).Value = 0.0002264463184709298
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.525608
$ws.Cells.Item(12, 8).Value = 1.576824
$ws.Cells.Item(12, 9).Value = 0.006785328120013172
$ws.Cells.Item(12, 10).Value = 0.006785328120013173
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.525608
$ws.Cells.Item(12, 14).Value = 1.576824
$ws.Cells.Item(12, 15).Value = 0.006785328120013172
$ws.Cells.Item(12, 16).Value = 0.006785328120013173
$ws.Cells.Item(12, 17).Value = 0.276263769664
$ws.Cells.Item(12, 18).Value = 2.486373926976
$ws.Cells.Item(12, 19).Value = [double]"4.604067769624149E-05"
$ws.Cells.Item(12, 20).Value = [double]"4.60406776962415E-05"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.525608
$ws.Cells.Item(13, 8).Value = 1.576824
$ws.Cells.Item(13, 9).Value = 0.006785328120013172
$ws.Cells.Item(13, 10).Value = 0.006785328120013173
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 34.11724066666667
$ws.Cells.Item(13, 14).Value = 102.351722
$ws.Cells.Item(13, 15).Value = 0.4404359759988248
$ws.Cells.Item(13, 16).Value = 0.4404359759988249
$ws.Cells.Item(13, 17).Value = 17.93229463232533
$ws.Cells.Item(13, 18).Value = 161.390651690928
$ws.Cells.Item(13, 19).Value = 0.002988502613010272
$ws.Cells.Item(13, 20).Value = 0.002988502613010274
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 34.11724066666667
$ws.Cells.Item(14, 8).Value = 102.351722
$ws.Cells.Item(14, 9).Value = 0.4404359759988248
$ws.Cells.Item(14, 10).Value = 0.4404359759988249
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 40.23443200000001
$ws.Cells.Item(14, 14).Value = 120.703296
$ws.Cells.Item(14, 15).Value = 0.5194057602668869
$ws.Cells.Item(14, 16).Value = 0.5194057602668869
$ws.Cells.Item(14, 17).Value = 1372.687799630635
$ws.Cells.Item(14, 18).Value = 12354.19019667571
$ws.Cells.Item(14, 19).Value = 0.2287649829625579
$ws.Cells.Item(14, 20).Value = 0.228764982962558
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 34.11724066666667
$ws.Cells.Item(15, 8).Value = 102.351722
$ws.Cells.Item(15, 9).Value = 0.4404359759988248
$ws.Cells.Item(15, 10).Value = 0.4404359759988249
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 2.585148666666667
$ws.Cells.Item(15, 14).Value = 7.755446
$ws.Cells.Item(15, 15).Value = 0.03337293561427507
$ws.Cells.Item(15, 16).Value = 0.03337293561427508
$ws.Cells.Item(15, 17).Value = 88.19813921977911
$ws.Cells.Item(15, 18).Value = 793.7832529780121
$ws.Cells.Item(15, 19).Value = 0.01469864146921918
$ws.Cells.Item(15, 20).Value = 0.01469864146921919
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 34.11724066666667
$ws.Cells.Item(16, 8).Value = 102.351722
$ws.Cells.Item(16, 9).Value = 0.4404359759988248
$ws.Cells.Item(16, 10).Value = 0.4404359759988249
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.525608
$ws.Cells.Item(16, 14).Value = 1.576824
$ws.Cells.Item(16, 15).Value = 0.006785328120013172
$ws.Cells.Item(16, 16).Value = 0.006785328120013173
$ws.Cells.Item(16, 17).Value = 17.93229463232533
$ws.Cells.Item(16, 18).Value = 161.390651690928
$ws.Cells.Item(16, 19).Value = 0.002988502613010272
$ws.Cells.Item(16, 20).Value = 0.002988502613010274
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 34.11724066666667
$ws.Cells.Item(17, 8).Value = 102.351722
$ws.Cells.Item(17, 9).Value = 0.4404359759988248
$ws.Cells.Item(17, 10).Value = 0.4404359759988249
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 34.11724066666667
$ws.Cells.Item(17, 14).Value = 102.351722
$ws.Cells.Item(17, 15).Value = 0.4404359759988248
$ws.Cells.Item(17, 16).Value = 0.4404359759988249
$ws.Cells.Item(17, 17).Value = 1163.986110707254
$ws.Cells.Item(17, 18).Value = 10475.87499636529
$ws.Cells.Item(17, 19).Value = 0.1939838489540374
$ws.Cells.Item(17, 20).Value = 0.1939838489540374
